$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.143.94'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '1.835.98'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  -2.17%  '
$ws.Range("D6").Value = "'0.6652"
$ws.Range("E6").Value = '  -4.43%  '
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'0.2957"
$ws.Range("E8").Value = '  -3.89%  '
$ws.Range("D9").Value = "'0.07364"
$ws.Range("E9").Value = '  -4.34%  '
$ws.Range("E10").Value = '  -3.69%  '
$ws.Range("D11").Value = "'0.07689"
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("D12").Value = '1.841.85'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = "'5.021"
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("D14").Value = "'0.6757"
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("D15").Value = "'86.22"
$ws.Range("E15").Value = '  -5.43%  '
$ws.Range("D16").Value = "'6.192"
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").Value = '29.060.37'
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("D18").Value = "'0.000008246"
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("D19").Value = "'228.81"
$ws.Range("E19").Value = '  -4.08%  '
$ws.Range("E20").Value = '  -1.84%  '
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = "'7.303"
$ws.Range("E22").Value = '  -4.14%  '
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = "'161.17"
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("E25").Value = '  -5.20%  '
$ws.Range("D26").Value = "'8.683"
$ws.Range("E26").Value = '  -2.44%  '
$ws.Range("D27").Value = "'18.04"
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("D28").Value = "'1.499"
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("D29").Value = "'4.235"
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").Value = "'4.105"
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("D31").Value = "'1.207"
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").Value = "'0.05308"
$ws.Range("E32").Value = '  +3.91%  '
$ws.Range("E33").Value = '  -1.12%  '
$ws.Range("D34").Value = "'0.7471"
$ws.Range("E34").Value = '  -3.83%  '
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("D36").Value = "'2.678"
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("D37").Value = '1.315.37'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = "'0.01806"
$ws.Range("E38").Value = '  -3.82%  '
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("D40").Value = "'0.9234"
$ws.Range("E40").Value = '  -3.35%  '
$ws.Range("E41").Value = '  +3.75%  '
$ws.Range("D42").Value = "'0.9987"
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("D43").Value = "'103.40"
$ws.Range("E43").Value = '  -2.48%  '
$ws.Range("D44").Value = '1.985.78'
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("D45").Value = "'0.5168"
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'63.82"
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = "'1.762"
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'9.288"
$ws.Range("E48").Value = '  -5.54%  '
$ws.Range("B49").Value = 'XinFinNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D49").Value = "'0.07440"
$ws.Range("E49").Value = '  +8.81%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.05931"
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = "'6.834"
$ws.Range("E51").Value = '  -1.82%  '
